$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Knight-Hennessey, Schwarzman" -> "Knight-Hennessy Schwarzman"
#    (drop the extra "e"/comma: "Hennessey, " becomes "Hennessy ")
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Knight-Hennessey, Schwarzman", $true, $false, $false, $false, $false, $true, 1, $false, "Knight-Hennessy Schwarzman", 2)

# ------------------------------------------------------------------
# 2. Split "Hennessy" into its own run, distinct from the rest of the
#    sentence (matches the source, which now carries this word in a
#    separate <w:r>).
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Hennessy")
$rng2.Font.Bold = $false

# ------------------------------------------------------------------
# 3. Section page size: make the (already-portrait) orientation
#    explicit on <w:pgSz>.
# ------------------------------------------------------------------
$d.PageSetup.Orientation = 0
